# Update the model-holdings workbook:
#  1. Update the "as of" date in the confidential disclosure note (A38) from
#     2021-06-10 to 2021-06-14.
#  2. Refresh the Weight (D) and Percent Change (E) values for rows 2-35.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The sheet is protected; unprotect it for the duration of the edit and
# re-apply protection afterwards.
$ws.Unprotect()

# --- 1. Update disclosure text date ---
$ws.Range("A38").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-06-14 for illustrative purposes only and are subject to change."

# --- 2. Update Weight / Percent Change values ---
$values = @{
    2  = @(0.08634016649012284, 0.02457793482528481)
    3  = @(0.07601872803857147, 0.007755244484082269)
    4  = @(0.04932757368671375, 0.01106718895193359)
    5  = @(0.05495521161896544, -0.01147060874416017)
    6  = @(0.04848172568909197, -0.003745318352059823)
    7  = @(0.04488654008079992, -0.01696924324661542)
    8  = @(0.03581773486357292, 0.007698954818533554)
    9  = @(0.03890045950519292, 0.002485451018428675)
    10 = @(0.03526577446087323, -0.001349911190053343)
    11 = @(0.03518603865904131, 0.003191836939857806)
    12 = @(0.03417615009314055, 0.01663346012195865)
    13 = @(0.0333157108859414, -0.01194457716196851)
    14 = @(0.03206219333588507, -0.004665830035074259)
    15 = @(0.03232174000321596, 0.00279069767441853)
    16 = @(0.03137198490708863, 0.00451009132934943)
    17 = @(0.02981470490310331, 0.0003707548568885333)
    18 = @(0.02696676608831926, 0.01278919149354807)
    19 = @(0.02431160810310296, 0.007032348804500765)
    20 = @(0.02196905572913527, -0.002616431187859702)
    21 = @(0.02290746210670242, -0.001608492842206877)
    22 = @(0.02317614670326358, -0.01086182336182351)
    23 = @(0.0221683215407374, 0.005877268798617052)
    24 = @(0.01987823856686613, -0.007844474761255049)
    25 = @(0.02199293225389825, -0.003377563329312383)
    26 = @(0.02130316598296778, -0.01295143212951422)
    27 = @(0.02098805481047861, -0.0109549023187876)
    28 = @(0.01781910938499231, 0.005186061322900448)
    29 = @(0.01750149265126257, -0.004395937547369888)
    30 = @(0.01050876600077201, 0.01085538772247241)
    31 = @(0.008001583514704032, 0.01282004052311669)
    32 = @(0.007203783338518861, 0.02275098717188051)
    33 = @(0.007977412218030401, 0.02898791708236348)
    34 = @(0.007083663784927339, 0.0247596854063501)
    35 = @(0.9999999999999999, 0.002751401014044808)
}

foreach ($row in $values.Keys) {
    $pair = $values[$row]
    $ws.Cells.Item($row, 4).Value = $pair[0]
    $ws.Cells.Item($row, 5).Value = $pair[1]
}

# --- 3. Restore sheet protection ---
$ws.Protect()
